# Applies the crypto price/volume table refresh described by the commit
# "Updated symbol list on Wed Jan  4 09:44:47 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / Link columns are plain text - a direct assignment is fine.
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# Price / Volume columns look numeric ("255.20", "3.69%", "1,671.99%", ...)
# but must stay plain text (leading/trailing zeros, "%" sign, "--" placeholders,
# etc. all need to be preserved verbatim). Force a text number format before
# assigning the value, then drop back to the default "Normal" style so the
# saved file does not end up with stray custom formatting on these cells.
function Set-NumericTextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-NumericTextCell $ws 'D2' '255.20'
Set-NumericTextCell $ws 'E2' '3.69%'

# Row 3
Set-NumericTextCell $ws 'D3' '28.11'
Set-NumericTextCell $ws 'E3' '-5.78%'

# Row 4
Set-NumericTextCell $ws 'D4' '5.249'
Set-NumericTextCell $ws 'E4' '1.91%'

# Row 5
Set-NumericTextCell $ws 'D5' '0.05852'
Set-NumericTextCell $ws 'E5' '1.71%'

# Row 6
Set-NumericTextCell $ws 'E6' '0.82%'

# Row 7
Set-TextCell $ws 'B7' 'GateToken'
Set-TextCell $ws 'C7' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-NumericTextCell $ws 'D7' '3.230'
Set-NumericTextCell $ws 'E7' '-0.32%'

# Row 8
Set-TextCell $ws 'B8' 'MXToken'
Set-TextCell $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-NumericTextCell $ws 'D8' '0.8682'
Set-NumericTextCell $ws 'E8' '2.29%'

# Row 9
Set-TextCell $ws 'B9' 'FTXToken'
Set-TextCell $ws 'C9' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-NumericTextCell $ws 'D9' '1.031'
Set-NumericTextCell $ws 'E9' '20.76%'

# Row 10
Set-TextCell $ws 'B10' 'WazirX'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-NumericTextCell $ws 'D10' '0.1412'
Set-NumericTextCell $ws 'E10' '1.75%'

# Row 11
Set-TextCell $ws 'B11' 'MandalaExchangeToken'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-NumericTextCell $ws 'D11' '0.07133'
Set-NumericTextCell $ws 'E11' '0.66%'

# Row 12
Set-TextCell $ws 'B12' 'BitrueCoin'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-NumericTextCell $ws 'D12' '0.03185'
Set-NumericTextCell $ws 'E12' '-2.17%'

# Row 13
Set-TextCell $ws 'B13' 'BitMartToken'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-NumericTextCell $ws 'D13' '0.09226'
Set-NumericTextCell $ws 'E13' '-1.58%'

# Row 14
Set-TextCell $ws 'B14' 'BitForexToken'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-NumericTextCell $ws 'D14' '0.001552'
Set-NumericTextCell $ws 'E14' '1.64%'

# Row 15
Set-NumericTextCell $ws 'D15' '0.005832'
Set-NumericTextCell $ws 'E15' '-1.46%'

# Row 16
Set-NumericTextCell $ws 'E16' '-0.66%'

# Row 17
Set-TextCell $ws 'B17' 'BTSEToken'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-NumericTextCell $ws 'D17' '2.222'
Set-NumericTextCell $ws 'E17' '0.02%'

# Row 18
Set-TextCell $ws 'B18' 'One'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-NumericTextCell $ws 'D18' '0.01061'
Set-NumericTextCell $ws 'E18' '1,671.99%'

# Row 19
Set-NumericTextCell $ws 'D19' '0.3183'
Set-NumericTextCell $ws 'E19' '0.56%'

# Row 20
Set-NumericTextCell $ws 'D20' '0.03483'
Set-NumericTextCell $ws 'E20' '3.66%'

# Row 21
Set-NumericTextCell $ws 'E21' '0.08%'

# Row 22
Set-NumericTextCell $ws 'D22' '3.546'
Set-NumericTextCell $ws 'E22' '1.42%'

# Row 23
Set-NumericTextCell $ws 'E23' '0.87%'

# Row 24
Set-NumericTextCell $ws 'E24' '-4.39%'

# Row 25
Set-NumericTextCell $ws 'D25' '0.001222'
Set-NumericTextCell $ws 'E25' '-0.59%'

# Row 26
Set-NumericTextCell $ws 'D26' '0.004880'
Set-NumericTextCell $ws 'E26' '17.76%'

# Row 27
Set-NumericTextCell $ws 'E27' '0.07%'

# Row 28
Set-NumericTextCell $ws 'E28' '33.86%'

# Row 40
Set-NumericTextCell $ws 'E40' '1.87%'

# Row 41
Set-NumericTextCell $ws 'D41' '0.005782'
Set-NumericTextCell $ws 'E41' '0.59%'

# Row 42
Set-NumericTextCell $ws 'D42' '0.1100'
Set-NumericTextCell $ws 'E42' '2.82%'

# Row 43
Set-NumericTextCell $ws 'D43' '0.002343'
Set-NumericTextCell $ws 'E43' '1.93%'

# Row 44
Set-NumericTextCell $ws 'D44' '0.009694'
Set-NumericTextCell $ws 'E44' '-2.70%'

# Row 45
Set-NumericTextCell $ws 'D45' '0.00005238'
Set-NumericTextCell $ws 'E45' '-5.16%'

# Row 46
Set-NumericTextCell $ws 'E46' '0.07%'

# Row 47
Set-NumericTextCell $ws 'E47' '31.09%'

# Row 48
Set-NumericTextCell $ws 'E48' '-12.79%'

# Row 49
Set-NumericTextCell $ws 'E49' '0.07%'

# Row 50
Set-NumericTextCell $ws 'E50' '0.07%'
